$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last one and name it "Partida 7"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Partida 7"

# Header row (same shared strings T / V / A as the other "Partida" sheets)
$ws.Range("A1").Value = "T"
$ws.Range("B1").Value = "V"
$ws.Range("C1").Value = "A"

# Data rows
$data = @(
    @(10, -5, 15),
    @(10, -5, 15),
    @(-5, 15, 10),
    @(15, -5, 25),
    @(-5, 10, 15),
    @(-5, 25, -5),
    @(-10, -5, -10),
    @(-10, -5, -5),
    @(10, -5, 25),
    @(-5, 35, -5),
    @(35, -10, -10),
    @(-10, -5, -10),
    @(-5, 30, -5),
    @(30, -10, 30),
    @(-5, 25, 45),
    @(30, -5, -5),
    @(25, 25, -5),
    @(-15, -5, -5),
    @(25, 20, -5),
    @(25, -5, -5),
    @(20, 25, -5),
    @(25, 15, -5),
    @(20, 20, -5),
    @(10, 25, -10),
    @(10, 10, 10),
    @(-5, 10, 10),
    @(15, -5, 10),
    @(15, 10, -5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Match the recorded selection on the new sheet
$ws.Range("H18").Select() | Out-Null
